{"js": "// Replace the date and every \"AAA\u00d7B=\" multiplication prompt in the\n// document with its updated value, per the commit diff. Every value\n// in the document is unique, so a plain case-sensitive search/replace\n// for each old->new pair is unambiguous.\nconst replacements = [\n  [\"2024-02-20 Tuesday\", \"2024-02-21 Wednesday\"],\n  [\"445\u00d74=\", \"108\u00d78=\"],\n  [\"395\u00d76=\", \"809\u00d75=\"],\n  [\"239\u00d73=\", \"378\u00d77=\"],\n  [\"942\u00d75=\", \"530\u00d77=\"],\n  [\"433\u00d77=\", \"299\u00d77=\"],\n  [\"203\u00d76=\", \"421\u00d74=\"],\n  [\"681\u00d78=\", \"225\u00d76=\"],\n  [\"248\u00d77=\", \"322\u00d75=\"],\n  [\"998\u00d77=\", \"889\u00d76=\"],\n  [\"802\u00d76=\", \"837\u00d79=\"],\n  [\"776\u00d79=\", \"635\u00d76=\"],\n  [\"234\u00d76=\", \"336\u00d72=\"],\n  [\"177\u00d74=\", \"221\u00d77=\"],\n  [\"796\u00d75=\", \"958\u00d79=\"],\n  [\"990\u00d75=\", \"227\u00d78=\"],\n  [\"249\u00d75=\", \"689\u00d78=\"],\n  [\"120\u00d78=\", \"673\u00d76=\"],\n  [\"176\u00d77=\", \"522\u00d74=\"],\n  [\"867\u00d78=\", \"116\u00d78=\"],\n  [\"331\u00d74=\", \"345\u00d72=\"],\n  [\"851\u00d79=\", \"537\u00d77=\"],\n  [\"604\u00d77=\", \"113\u00d74=\"],\n  [\"224\u00d74=\", \"112\u00d72=\"],\n  [\"494\u00d74=\", \"679\u00d73=\"],\n  [\"460\u00d74=\", \"849\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and every \"AAA\u00d7B=\" multiplication prompt in the\n# document with its updated value, per the commit diff. Every value\n# in the document is unique, so a plain case-sensitive Find/Replace\n# for each old->new pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-20 Tuesday\", \"2024-02-21 Wednesday\"),\n    @(\"445\u00d74=\", \"108\u00d78=\"),\n    @(\"395\u00d76=\", \"809\u00d75=\"),\n    @(\"239\u00d73=\", \"378\u00d77=\"),\n    @(\"942\u00d75=\", \"530\u00d77=\"),\n    @(\"433\u00d77=\", \"299\u00d77=\"),\n    @(\"203\u00d76=\", \"421\u00d74=\"),\n    @(\"681\u00d78=\", \"225\u00d76=\"),\n    @(\"248\u00d77=\", \"322\u00d75=\"),\n    @(\"998\u00d77=\", \"889\u00d76=\"),\n    @(\"802\u00d76=\", \"837\u00d79=\"),\n    @(\"776\u00d79=\", \"635\u00d76=\"),\n    @(\"234\u00d76=\", \"336\u00d72=\"),\n    @(\"177\u00d74=\", \"221\u00d77=\"),\n    @(\"796\u00d75=\", \"958\u00d79=\"),\n    @(\"990\u00d75=\", \"227\u00d78=\"),\n    @(\"249\u00d75=\", \"689\u00d78=\"),\n    @(\"120\u00d78=\", \"673\u00d76=\"),\n    @(\"176\u00d77=\", \"522\u00d74=\"),\n    @(\"867\u00d78=\", \"116\u00d78=\"),\n    @(\"331\u00d74=\", \"345\u00d72=\"),\n    @(\"851\u00d79=\", \"537\u00d77=\"),\n    @(\"604\u00d77=\", \"113\u00d74=\"),\n    @(\"224\u00d74=\", \"112\u00d72=\"),\n    @(\"494\u00d74=\", \"679\u00d73=\"),\n    @(\"460\u00d74=\", \"849\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # MatchCase:=True, Forward:=True, Wrap:=wdFindContinue(1), Replace:=wdReplaceAll(2)\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
